# Addresses issue #33: add placeholder time-of-day (0 = midnight) for mass
# measurement rows that were missing a "B" (time) value, so Excel stops
# warning about them. Also restores the scroll/selection position on the
# "valvemap.csv" sheet to the top of the data (A3 / B3) instead of the
# bottom where it had been left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valvemap.csv")

# Reference cell that already carries the "populated time" style (font +
# h:mm number format) used throughout column B wherever a reading exists.
$styleSource = $ws.Range("B376")
$styleSource.Copy()

# Contiguous blocks of rows whose "B" cell exists but is blank (style-only,
# no <v>) -- fill them with placeholder 0 (12:00 AM) using the same style
# already used for populated time cells.
$blocksToFill = @(
    "B143:B172",
    "B210:B239",
    "B276:B305",
    "B377:B406",
    "B450:B479",
    "B515:B544",
    "B616:B645",
    "B682:B711",
    "B747:B776",
    "B814:B843",
    "B880:B909",
    "B946:B975",
    "B1018:B1047",
    "B1084:B1113",
    "B1150:B1179"
)

foreach ($block in $blocksToFill) {
    $dst = $ws.Range($block)
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $dst.Value = 0
}

$excel.CutCopyMode = 0

# Restore the view: unfreeze-pane scroll position back to the top of the
# data and the active selection to B3 (it had drifted down near row 1196).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B3").Select()
